$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Programme / Department rows -----------------------------------------
# Re-point the existing rows to the new order and fill in the Department
# column (previously blank) for each programme.
$ws.Range("A2").Value = "Tech Fest"
$ws.Range("B2").Value = "Department of Computer Science"

$ws.Range("A3").Value = "Culture Fest"
$ws.Range("B3").Value = "Department of Cultural Studies"

$ws.Range("A4").Value = "Commerce Fest"
$ws.Range("B4").Value = "Department of Commerce"

$ws.Range("A5").Value = "Hackathon"
$ws.Range("B5").Value = "Department of Computer Science"

# New programme row
$ws.Range("A6").Value = "Movie Fest"

# --- Extend the sheet's used range down to row 101 ------------------------
# Touch each cell's formatting (without actually changing it) so Excel
# materialises the (blank) cells and the sheet's dimension grows to match.
$ws.Range("A7:A101").Font.Bold = $false

# --- Department dropdown list ----------------------------------------------
# Add the new "Department of Media Communications" option to the list
# validation applied to B2:B100.
$rng = $ws.Range("B2:B100")
$v = $rng.Validation
$v.Delete()
$v.Add(3, 1, 1, '"Department of Computer Science, Department of Commerce, Department of Business Studies, Department of Cultural Studies, Department of Media Communications"')
$v.IgnoreBlank = $false

Write-Host "done"
